$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: split "seismoacoustic ... input - categorical" run into 3 runs
#           and add a new "shift - input - binary" paragraph right after it.
# ---------------------------------------------------------------------------

$enDash = [string][char]0x2013

# Locate the unique "input - categorical" text (plain hyphen identifies the
# seismoacoustic bullet, since the other similar bullets use an en-dash there).
$target = $d.Content.Duplicate
$target.Start = 0
$target.End = $d.Content.End
$target.Find.ClearFormatting()
$target.Find.Execute("input - categorical") | Out-Null
if (-not $target.Find.Found) { throw "Could not find 'input - categorical' text" }

# Expand left by 3 characters to capture the leading " <en-dash> " that is
# part of the very same run.
$fullRun = $d.Range($target.Start - 3, $target.End)
$startPos = $fullRun.Start

# Wipe the run's text, then rebuild it as three separate runs via sequential
# InsertAfter calls (each InsertAfter on a fresh collapsed Range creates a
# new run boundary).
$fullRun.Text = ""

$r1 = $d.Range($startPos, $startPos)
$r1.InsertAfter(" " + $enDash + " input ")
$pos = $startPos + 9

$r2 = $d.Range($pos, $pos)
$r2.InsertAfter($enDash)
$pos = $pos + 1

$r3 = $d.Range($pos, $pos)
$r3.InsertAfter(" categorical")
$pos = $pos + 12

# Paragraph that now holds "seismoacoustic ... categorical"
$seismoPara = $d.Range($pos, $pos).Paragraphs(1)

# Insert a brand-new paragraph right after it containing the "shift" bullet.
$seismoPara.Range.InsertParagraphAfter()
$shiftInsertPos = $seismoPara.Range.End
$shiftRng = $d.Range($shiftInsertPos, $shiftInsertPos)
$shiftRng.InsertAfter("shift " + $enDash + " input " + $enDash + " binary ")

# ---------------------------------------------------------------------------
# Change 2: add a new "Main Effects" narrative paragraph, and move the
#           "_GoBack" bookmark paragraph from before the heading to after
#           the new narrative paragraph.
# ---------------------------------------------------------------------------

# Locate the "Main Effects" Heading3 paragraph.
$mainEffectsPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Main Effects" + [char]13) { $mainEffectsPara = $p; break }
}
if ($mainEffectsPara -eq $null) { throw "Could not find 'Main Effects' paragraph" }

# Insert a new blank paragraph directly above the heading (this is the blank
# line that ends up between the "skewed" paragraph and the heading once the
# bookmark paragraph is relocated below).
$mainEffectsPara.Range.InsertParagraphBefore()
# Re-fetch the heading paragraph (InsertParagraphBefore can shift references)
# and make sure the freshly inserted blank paragraph has no stray style.
$mainEffectsPara2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Main Effects" + [char]13) { $mainEffectsPara2 = $p; break }
}
$blankBefore = $mainEffectsPara2.Previous()
$blankBefore.Range.Style = "Normal"

# Relocate the "_GoBack" bookmark: delete it from its current paragraph and
# remove that now-empty paragraph entirely.
$bm = $d.Bookmarks.Item("_GoBack")
$bmRange = $bm.Range
$bmParaStart = $bmRange.Start
$bm.Delete()
$emptyBmPara = $d.Range($bmParaStart, $bmParaStart).Paragraphs(1)
$emptyBmPara.Range.Delete()

# Locate the "The 'nbumps' " paragraph that follows "Main Effects" and
# replace its content with the new narrative text (including proofErr
# markers), using InsertXML for precise run-level control.
$nbumpsPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -eq "The " + [char]0x2018 + "nbumps" + [char]0x2019 + " " + [char]13) { $nbumpsPara = $p; break }
}
if ($nbumpsPara -eq $null) { throw "Could not find the 'nbumps' paragraph" }

$nbumpsRng = $nbumpsPara.Range.Duplicate
$nbumpsRng.MoveEnd(1, -1)
$nbumpsRng.Text = ""

$insertRng = $nbumpsPara.Range.Duplicate
$insertRng.MoveEnd(1, -1)

$mainEffectsXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t xml:space="preserve">The </w:t></w:r>
<w:r><w:t xml:space="preserve">main effects for this study </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>are considered to be</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t xml:space="preserve"> the all numeric variables, plus </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>ghazard</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve">, </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>seismoacoustic</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> and shift.  </w:t></w:r>
<w:r><w:t>This stands to reason, since numerical energy readings and shift activity type all seem like they would impact the number of hazardous seismic events in the next shift.  The nbumps class of variables are left out for more advanced models, since the resonance and frequency ranges could have a multitude of confounding variables that we, without significant mining expertise, would miss.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$insertRng.InsertXML($mainEffectsXml)

# Insert a new blank paragraph after the narrative paragraph, then re-add the
# "_GoBack" bookmark in its own paragraph right after that.
$narrativePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*without significant mining expertise, would miss.*") { $narrativePara = $p; break }
}
if ($narrativePara -eq $null) { throw "Could not find the new narrative paragraph" }

$narrativePara.Range.InsertParagraphAfter()
$narrativePara2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*without significant mining expertise, would miss.*") { $narrativePara2 = $p; break }
}
$blankAfter = $narrativePara2.Next()
$blankAfter.Range.Style = "Normal"

$blankAfter.Range.InsertParagraphAfter()
$bmTargetPara = $blankAfter.Next()
$bmInsertRng = $bmTargetPara.Range.Duplicate
$bmInsertRng.Collapse(1)
$bmInsertRng.Style = "Normal"
$d.Bookmarks.Add("_GoBack", $bmInsertRng) | Out-Null

Write-Output "Edit complete."
